# Apply cell value updates from the cryptos list refresh (GitHub Actions run).
# Values are written via NumberFormat "@" (Text) so numeric-looking strings
# (prices, percentages) stay literal text, then the format/style is reset
# back to the default "Normal" style to avoid introducing new cell styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "70.913.61"
Set-TextValue "E2" "  -0.23%  "
Set-TextValue "D3" "3.844.68"
Set-TextValue "E3" "  +0.75%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "699.82"
Set-TextValue "E5" "  -1.21%  "
Set-TextValue "D6" "172.01"
Set-TextValue "E6" "  -0.87%  "
Set-TextValue "D7" "3.840.61"
Set-TextValue "E7" "  +0.66%  "
Set-TextValue "E9" "  -0.50%  "
Set-TextValue "E10" "  -1.45%  "
Set-TextValue "E11" "  -1.60%  "
Set-TextValue "E12" "  -0.87%  "
Set-TextValue "E13" "  -1.60%  "
Set-TextValue "D14" "36.23"
Set-TextValue "E14" "  -0.30%  "
Set-TextValue "D15" "4.495.40"
Set-TextValue "E15" "  +0.83%  "
Set-TextValue "D16" "3.991.34"
Set-TextValue "E16" "  +4.78%  "
Set-TextValue "D17" "70.955.86"
Set-TextValue "E17" "  -0.26%  "
Set-TextValue "B18" "TRON"
Set-TextValue "C18" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.115"
Set-TextValue "E18" "  +0.60%  "
Set-TextValue "B19" "Polkadot"
Set-TextValue "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "7.16"
Set-TextValue "E19" "  -1.45%  "
Set-TextValue "E20" "  -3.37%  "
Set-TextValue "E21" "  -4.77%  "
Set-TextValue "D22" "493.43"
Set-TextValue "E22" "  +1.88%  "
Set-TextValue "E23" "  -0.28%  "
Set-TextValue "D24" "84.55"
Set-TextValue "E24" "  +0.67%  "
Set-TextValue "D25" "0.0000147"
Set-TextValue "E25" "  +0.36%  "
Set-TextValue "B26" "InternetComputer(DFINITY)"
Set-TextValue "C26" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D26" "12.15"
Set-TextValue "E26" "  -2.80%  "
Set-TextValue "B27" "RenderToken"
Set-TextValue "C27" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D27" "10.54"
Set-TextValue "E27" "  -0.82%  "
Set-TextValue "E28" "  -3.62%  "
Set-TextValue "E29" "  +0.05%  "
Set-TextValue "D30" "3.15"
Set-TextValue "E30" "  +1.59%  "
Set-TextValue "E31" "  -0.78%  "
Set-TextValue "E32" "  -2.17%  "
Set-TextValue "D33" "29.42"
Set-TextValue "E33" "  -0.71%  "
Set-TextValue "E34" "  +1.06%  "
Set-TextValue "D35" "3.803.10"
Set-TextValue "E35" "  +0.97%  "
Set-TextValue "D36" "9.13"
Set-TextValue "E36" "  -1.49%  "
Set-TextValue "D37" "1.00"
Set-TextValue "E37" "  -0.05%  "
Set-TextValue "E38" "  -1.16%  "
Set-TextValue "D39" "2.39"
Set-TextValue "E39" "  +6.01%  "
Set-TextValue "E40" "  +6.49%  "
Set-TextValue "E41" "  -0.16%  "
Set-TextValue "E42" "  -4.97%  "
Set-TextValue "E44" "  -0.03%  "
Set-TextValue "E45" "  -9.01%  "
Set-TextValue "E46" "  +1.52%  "
Set-TextValue "D47" "48.82"
Set-TextValue "E47" "  -1.23%  "
Set-TextValue "E48" "  -1.08%  "
Set-TextValue "B49" "Arweave"
Set-TextValue "C49" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D49" "43.36"
Set-TextValue "E49" "  -4.63%  "
Set-TextValue "B50" "Cosmos"
Set-TextValue "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "8.62"
Set-TextValue "E50" "  +0.37%  "
Set-TextValue "B51" "Bittensor"
Set-TextValue "C51" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D51" "407.67"
Set-TextValue "E51" "  +0.97%  "
